$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "Test ID"
$ws.Range("B1").Value = "Collection Date"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"
$ws.Range("E1").Value = "Name"
$ws.Range("F1").Value = "Area (ha)"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "Age"
$ws.Range("I1").Value = "Address"
$ws.Range("J1").Value = "Mobile No."
$ws.Range("K1").Value = "Soil pH"
$ws.Range("L1").Value = "Nitrogen"
$ws.Range("M1").Value = "Phosphorus"
$ws.Range("N1").Value = "Potassium"
$ws.Range("O1").Value = "Electrical Conductivity"
$ws.Range("P1").Value = "Temperature"
$ws.Range("Q1").Value = "Moisture"
$ws.Range("R1").Value = "Humidity"
$ws.Range("S1").Value = "Soil Health Score"
$ws.Range("T1").Value = "Recommendations"
$ws.Range("U1").Value = "Fertilizer Recommendation"

# --- Update data row (row 2) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2659"
$ws.Range("B2").Value = "25-04-2024"
$ws.Range("C2").Value = 265.256598
$ws.Range("D2").Value = 352.356565
$ws.Range("E2").Value = "Lalchawimawii"
$ws.Range("F2").Value = 2.5
$ws.Range("G2").Value = "Female"
$ws.Range("H2").Value = 33
$ws.Range("I2").Value = "Zarkawt, Aizawl"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "9865323265"
$ws.Range("K2").Value = 5.5
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 200
$ws.Range("O2").Value = 3.5
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 30
$ws.Range("R2").Value = 40
$ws.Range("S2").Value = 0.5535232709427351
$ws.Range("T2").Value = "Grow maize, soybean, groundnut, cotton, and incorporate legumes into the cropping system."
$ws.Range("U2").Value = "No specific fertilizer recommendation available for the given soil data. Please consult with local agriculture experts for personalized recommendations."
